# Applies the "NOX" abatement-data update to Othdata_dors_3.xlsx
#
# Sheets: 1=Y, 2=Q2P, 3=M, 4=M_sets
#
# Changes:
#  - M_sets (sheet4): append a new "NOX" value to the z-set list (A3), becomes
#    the active sheet/tab and selection.
#  - M (sheet3): fill in A3/B3 ("NOX", 10) next to the existing row-3 data,
#    and append 4 new NOX rows (electricity/oil/inp3/K) below the data.
#  - Y (sheet1): no longer the active tab; selection moves to E12.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$wsY     = $sheets.Item("Y")
$wsM     = $sheets.Item("M")
$wsMsets = $sheets.Item("M_sets")

# ---------------------------------------------------------------------------
# M_sets: add "NOX" as a new set member in column A, row 3.
# ---------------------------------------------------------------------------
$wsMsets.Range("A3").Value = "NOX"

# ---------------------------------------------------------------------------
# M: fill A3/B3 ("NOX", 10) alongside the existing row-3 data (no row shift),
# then append the four new NOX data rows below the existing data (rows 6-9).
# ---------------------------------------------------------------------------
$wsM.Range("A3").Value = "NOX"
$wsM.Range("B3").Value = 10

$wsM.Range("C6").Value = "NOX"
$wsM.Range("D6").Value = "electricity"
$wsM.Range("E6").Value = 2

$wsM.Range("C7").Value = "NOX"
$wsM.Range("D7").Value = "oil"
$wsM.Range("E7").Value = 0

$wsM.Range("C8").Value = "NOX"
$wsM.Range("D8").Value = "inp3"
$wsM.Range("E8").Value = 0

$wsM.Range("C9").Value = "NOX"
$wsM.Range("D9").Value = "K"
$wsM.Range("E9").Value = 0

# ---------------------------------------------------------------------------
# Selections / active sheet.
# ---------------------------------------------------------------------------
$wsY.Range("E12").Select()
$wsM.Range("C10").Select()
$wsMsets.Range("A4").Select()

$wsMsets.Activate()
